$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.402.41"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "3.522.90"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'612.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'151.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").Value = "3.520.90"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "'0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "'32.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "4.117.84"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "3.516.34"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "67.398.28"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'15.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "'445.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").Value = "'9.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'0.625"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").Value = "'77.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "'0.0000130"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.21%  "
$ws.Range("D26").Value = "3.662.10"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "'10.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'8.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'2.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -6.91%  "
$ws.Range("D33").Value = "'0.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("D34").Value = "'25.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").Value = "3.513.08"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D40").Value = "'177.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'2.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "'28.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").Value = "'44.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "'2.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D51").Value = "'0.997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
